# Append a new data row (row 73) to each of the 4 worksheets, matching the
# structure/style of the existing rows (date in col A with the custom
# "YYYY-MM-DD HH:MM:SS" format, hex-byte strings in cols B-E, numbers in F-I).

$wb = $excel.ActiveWorkbook

$dateValue = 45859.46223379629
$dateFormat = "YYYY-MM-DD HH:MM:SS"

# Sheet 1: MID_LFT_#1
$ws = $wb.Worksheets.Item(1)
$r = 73
$ws.Cells.Item($r, 1).Value = $dateValue
$ws.Cells.Item($r, 1).NumberFormat = $dateFormat
$ws.Cells.Item($r, 2).Value = "0x01,0x90"
$ws.Cells.Item($r, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"
$ws.Cells.Item($r, 4).Value = "0x01,0x44"
$ws.Cells.Item($r, 5).Value = "0x07"
$ws.Cells.Item($r, 6).Value = 400
$ws.Cells.Item($r, 7).Value = [double]"5.68631262647113e+23"
$ws.Cells.Item($r, 8).Value = 324
$ws.Cells.Item($r, 9).Value = 7

# Sheet 2: MID_LFT_#2
$ws = $wb.Worksheets.Item(2)
$r = 73
$ws.Cells.Item($r, 1).Value = $dateValue
$ws.Cells.Item($r, 1).NumberFormat = $dateFormat
$ws.Cells.Item($r, 2).Value = "0x01,0x7c"
$ws.Cells.Item($r, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item($r, 4).Value = "0x01,0x40"
$ws.Cells.Item($r, 5).Value = "0x19"
$ws.Cells.Item($r, 6).Value = 380
$ws.Cells.Item($r, 7).Value = [double]"5.68432987514711e+23"
$ws.Cells.Item($r, 8).Value = 320
$ws.Cells.Item($r, 9).Value = 25

# Sheet 3: MID_PLT_#1
$ws = $wb.Worksheets.Item(3)
$r = 73
$ws.Cells.Item($r, 1).Value = $dateValue
$ws.Cells.Item($r, 1).NumberFormat = $dateFormat
$ws.Cells.Item($r, 2).Value = "0x00,0x6e"
$ws.Cells.Item($r, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
$ws.Cells.Item($r, 4).Value = "0x00,0x63"
$ws.Cells.Item($r, 5).Value = "0x15"
$ws.Cells.Item($r, 6).Value = 110
$ws.Cells.Item($r, 7).Value = [double]"5.68631262647113e+23"
$ws.Cells.Item($r, 8).Value = 99
$ws.Cells.Item($r, 9).Value = 15

# Sheet 4: MID_PLT_#2
$ws = $wb.Worksheets.Item(4)
$r = 73
$ws.Cells.Item($r, 1).Value = $dateValue
$ws.Cells.Item($r, 1).NumberFormat = $dateFormat
$ws.Cells.Item($r, 2).Value = "0x00,0x82"
$ws.Cells.Item($r, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
$ws.Cells.Item($r, 4).Value = "0x00,0x78"
$ws.Cells.Item($r, 5).Value = "0x9"
$ws.Cells.Item($r, 6).Value = 130
$ws.Cells.Item($r, 7).Value = [double]"5.68631262647113e+23"
$ws.Cells.Item($r, 8).Value = 120
$ws.Cells.Item($r, 9).Value = 9
